$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (28) down to the new row (29)
$ws.Range("A28:B28").Copy()
$ws.Range("A29:B29").PasteSpecial(-4122)

# Set the new values for the appended row
$ws.Range("A29").Value = "20-10-2025"
$ws.Range("B29").Value = "The price of gold in India today is ₹13,069 per gram for 24 karat gold, ₹11,980 per gram for 22 karat gold and ₹9,802 per gram for 18 karat gold (also called 999 gold)."
